# "add tutorial project to folder" - To Touch the Moon.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Character Creation
# Insert 4 new blank columns before column H (old H:Q -> new L:U),
# then fill in a new "Effect" column (C) describing each stat, add a
# new "luck" stat row, and a small Hp/Mp/Ep legend block.
# ---------------------------------------------------------------
$wsCC = $wb.Worksheets.Item("Character Creation")

$wsCC.Range("H1:K1").EntireColumn.Insert() | Out-Null

$wsCC.Range("C6").Value = "melee damage, main carrying capacity"
$wsCC.Range("C7").Value = "helps to learn magic, increases mp regen"
$wsCC.Range("C8").Value = "increase skill die, high wisdom can grants some minor luck"
$wsCC.Range("C9").Value = "shop discount"
$wsCC.Range("C10").Value = "missile damage, minor increase in evasion"
$wsCC.Range("C11").Value = "main hp increase, minor carrying capacity"
$wsCC.Range("C12").Value = "defence against lots of monster abilities, minor hp and mp"
$wsCC.Range("C14").Value = "increase Mp, protect against some special abilities"

$wsCC.Range("A15").Value = "luck"
$wsCC.Range("C15").Value = "lots of minor things, helps crit chance evade traps,"

$wsCC.Range("I18").Value = "Hp"
$wsCC.Range("J18").Value = "Hit points"
$wsCC.Range("I19").Value = "Mp"
$wsCC.Range("J19").Value = "Mana points"
$wsCC.Range("I20").Value = "Ep"
$wsCC.Range("J20").Value = "Energy points"

$wsCC.Activate()
$wsCC.Range("L3:U13").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: procedurally generated
# Draw a small diamond-shaped room outline ("#" tiles).
# ---------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("procedurally generated")

$wsProc.Range("K13").Value = "#"
$wsProc.Range("L13").Value = "#"
$wsProc.Range("M13").Value = "#"

$wsProc.Range("J14").Value = "#"
$wsProc.Range("K14").Value = "#"
$wsProc.Range("M14").Value = "#"
$wsProc.Range("N14").Value = "#"

$wsProc.Range("I15").Value = "#"
$wsProc.Range("J15").Value = "#"
$wsProc.Range("N15").Value = "#"

$wsProc.Range("I16").Value = "#"
$wsProc.Range("O16").Value = "#"

$wsProc.Range("J17").Value = "#"
$wsProc.Range("N17").Value = "#"
$wsProc.Range("O17").Value = "#"

$wsProc.Range("J18").Value = "#"
$wsProc.Range("K18").Value = "#"
$wsProc.Range("M18").Value = "#"
$wsProc.Range("N18").Value = "#"

$wsProc.Range("K19").Value = "#"
$wsProc.Range("L19").Value = "#"
$wsProc.Range("M19").Value = "#"

$wsProc.Activate()
$wsProc.Range("F9").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: Lore
# ---------------------------------------------------------------
$wsLore = $wb.Worksheets.Item("Lore")
$wsLore.Range("S1").EntireColumn.ColumnWidth = 10.75
$wsLore.Activate()
$wsLore.Range("X2").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: Controls
# ---------------------------------------------------------------
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Activate()
$wsControls.Range("E9").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: Quests
# Insert 2 rows above the "Lord of the dead Empire" block, add the
# reward note, and annotate who Clive/Matilda are.
# ---------------------------------------------------------------
$wsQuests = $wb.Worksheets.Item("Quests")

$wsQuests.Range("A9:A10").EntireRow.Insert() | Out-Null
$wsQuests.Range("A8").Value = "Reward: some of there old equipment"
$wsQuests.Range("E3").Value = "Clive is a Calvary / horse rider"
$wsQuests.Range("H3").Value = "Matilda is a Valkryie"

$wsQuests.Activate()
$wsQuests.Range("J15").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: Names
# New phrase + make this the active tab (matches the saved workbook
# view's activeTab moving from "procedurally generated" to "Names").
# ---------------------------------------------------------------
$wsNames = $wb.Worksheets.Item("Names")
$wsNames.Range("T4").Value = "It shall be engraved on your soul"

$wsNames.Activate()
$wsNames.Range("T4").Select() | Out-Null
